$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("Q3").Value = "Area"
$ws.Range("R3").Value = 257.69726249999991

Write-Output "done"
